# Roboflow Annotation Report 7/30/2025
# Append this week's weekly-progress row to the "Avances Etiquetado
# Roboflow" tracking table (Table1 on Sheet1), which currently ends at
# row 68 (date 29/7/2031, the last logged entry).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

$lastRow = $tbl.Range.Row + $tbl.Range.Rows.Count - 1   # 68
$newRowIndex = $lastRow + 1                              # 69

# Clone the formatting (cell styles + row height) of the last data row
# down into the new row before writing values into it, so the new row
# keeps the same borders/number formats as the rest of the table.
$srcRange = $ws.Range("D$lastRow`:J$lastRow")
$dstRange = $ws.Range("D$newRowIndex`:J$newRowIndex")
$srcRange.Copy($dstRange)
$ws.Rows.Item($newRowIndex).RowHeight = $ws.Rows.Item($lastRow).RowHeight

# New weekly entry.
$ws.Range("D$newRowIndex").Value = "30/7/2032"
$ws.Range("E$newRowIndex").Value = 329
$ws.Range("F$newRowIndex").Value = 1001
$ws.Range("G$newRowIndex").Value = 0
$ws.Range("H$newRowIndex").Value = 0
$ws.Range("I$newRowIndex").Value = 1012
$ws.Range("J$newRowIndex").Value = "N/A"

# Grow the table (ref/autoFilter) and the sheet's used range follow.
$tbl.Resize($ws.Range("D4:J$newRowIndex"))
